# Negate the "Block" (column E) values for every data row on the sheet.
# Blank cells and zero values are left untouched since negation is a no-op
# for them (and keeps blank cells blank rather than turning them into 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A ("Command") is populated on every data row, so anchor the
# last-row lookup there (xlUp = -4162) rather than on column E, which has
# several blank/0 cells scattered throughout (including possibly at the
# very bottom of the used range).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "" -and $val -ne 0) {
        $cell.Value = -1 * $val
    }
}
